# Deploying to gh-pages: add 2022 column (S) to the statistics table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 1.9210869108320343

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 1.020872301352429

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 2.8415499553180767

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Font.Bold = $true
$ws.Range("S7").Value = 1.5924017665043597

$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 2.5011433798307796

$ws.Range("R9").Copy()
$ws.Range("S9").PasteSpecial(-4122)
$ws.Range("S9").Value = 0.70098698968147144

$ws.Range("R10").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("S10").Font.Bold = $true
$ws.Range("S10").Value = 2.2312343573160249

$ws.Range("R11").Copy()
$ws.Range("S11").PasteSpecial(-4122)
$ws.Range("S11").Value = 2.4764236727529938

$ws.Range("R12").Copy()
$ws.Range("S12").PasteSpecial(-4122)
$ws.Range("S12").Value = 1.9888745417939038

$ws.Range("R13").Copy()
$ws.Range("S13").PasteSpecial(-4122)
$ws.Range("S13").Font.Bold = $true
$ws.Range("S13").Value = 1.3057776932131271

$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)
$ws.Range("S14").Value = 2.6056788910230639

$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial(-4122)
$ws.Range("S15").Value = 0

$ws.Range("R16").Copy()
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("S16").Font.Bold = $true
$ws.Range("S16").Value = 0.65058422463372112

$ws.Range("R17").Copy()
$ws.Range("S17").PasteSpecial(-4122)
$ws.Range("S17").Value = 0.65686622262510019

$ws.Range("R18").Copy()
$ws.Range("S18").PasteSpecial(-4122)
$ws.Range("S18").Value = 0.64442124527961442

$ws.Range("R19").Copy()
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("S19").Font.Bold = $true
$ws.Range("S19").Value = 2.5553368555544047

$ws.Range("R20").Copy()
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("S20").Value = 1.807815324711445

$ws.Range("R21").Copy()
$ws.Range("S21").PasteSpecial(-4122)
$ws.Range("S21").Value = 3.2928586128833093

$ws.Range("R22").Copy()
$ws.Range("S22").PasteSpecial(-4122)
$ws.Range("S22").Font.Bold = $true
$ws.Range("S22").Value = 1.8387963974300983

$ws.Range("R23").Copy()
$ws.Range("S23").PasteSpecial(-4122)
$ws.Range("S23").Value = 2.2260807622100529

$ws.Range("R24").Copy()
$ws.Range("S24").PasteSpecial(-4122)
$ws.Range("S24").Value = 1.4582467499325562

$ws.Range("R25").Copy()
$ws.Range("S25").PasteSpecial(-4122)
$ws.Range("S25").Font.Bold = $true
$ws.Range("S25").Value = 1.2245886088767601

$ws.Range("R26").Copy()
$ws.Range("S26").PasteSpecial(-4122)
$ws.Range("S26").Value = 1.3105423773238725

$ws.Range("R27").Copy()
$ws.Range("S27").PasteSpecial(-4122)
$ws.Range("S27").Value = 1.1375464261135158

$ws.Range("R28").Copy()
$ws.Range("S28").PasteSpecial(-4122)
$ws.Range("S28").Font.Bold = $true
$ws.Range("S28").Value = 2.4791112740241377

$ws.Range("R29").Copy()
$ws.Range("S29").PasteSpecial(-4122)
$ws.Range("S29").Value = 2.4279584268771761

$ws.Range("R30").Copy()
$ws.Range("S30").PasteSpecial(-4122)
$ws.Range("S30").Value = 2.5408788313520994

$ws.Range("R31").Copy()
$ws.Range("S31").PasteSpecial(-4122)
$ws.Range("S31").Font.Bold = $true
$ws.Range("S31").Value = 1.1238322680339958

$ws.Range("R32").Copy()
$ws.Range("S32").PasteSpecial(-4122)
$ws.Range("S32").Value = 0.57553956834532372

$ws.Range("R33").Copy()
$ws.Range("S33").PasteSpecial(-4122)
$ws.Range("S33").Value = 1.6467682173734046

$excel.CutCopyMode = $false
[void]$ws.Range("T3").Select()